$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

$ws1.Range("A3").Value = "Shubhi"
$ws1.Range("B3").Value = "Srivastava"
$ws1.Range("C3").Value = "001100"
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Utkarsh"
$ws1.Range("B4").Value = "Srivastava"
$ws1.Range("C4").Value = "001100"
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "Geetika"
$ws1.Range("B5").Value = "Srivastava"
$ws1.Range("C5").Value = "001100"
$ws1.Range("D5").Value = "Customer added successfully"

$ws2 = $wb.Worksheets.Add()
$ws2.Name = "OpenAccountTest"
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Sagrika Srivastava"
$ws2.Range("B2").Value = "Rupee"
